$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow the data refresh, then
# re-protect once all the edits are in place.
$ws.Unprotect()

# Update the confidential disclosure date text in A10 (2021-04-22 -> 2021-04-23)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."
# Re-fit the row height back to its natural size - writing the multi-line
# string otherwise leaves a stale explicit row height behind.
$ws.Rows.Item(10).AutoFit()

# Refresh Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2531157521263684
$ws.Range("E2").Value = 0.013348526778447

$ws.Range("D3").Value = 0.4902962108496888
$ws.Range("E3").Value = 0.007524858908895427

$ws.Range("D4").Value = 0.1011379862992283
$ws.Range("E4").Value = 0.01564561734213044

$ws.Range("D5").Value = 0.09899579735792655
$ws.Range("E5").Value = 0.01262272089761596

$ws.Range("D6").Value = 0.05645425336678785
$ws.Range("E6").Value = 0.01593444115638531

$ws.Range("E7").Value = 0.01079966173629354

# Restore sheet protection to its original state.
$ws.Protect()
